# "major accuracy check update"
#
# 1. The sample id shared string "E7420" becomes "E7420L" (this string is
#    shown in column G, rows 2-49, all of which share the same text).
# 2. The scrolled/top-left visible cell of the sheet view moves from A16 to
#    A15 (best effort - not all view-state round trips through this host,
#    but setting it is harmless).
# 3. Column H (rows 2-49) held a volatile `=FALSE()` formula that evaluated
#    to boolean FALSE; it is replaced with the literal boolean value FALSE
#    (no formula), matching the custom "TRUE"/"FALSE" number format already
#    applied to those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the repeated sample-id text shown in G2:G49.
$ws.Range("G2:G49").Value = "E7420L"

# 2. Scroll the view so row 15 is the top visible row (was row 16).
try {
    $excel.ActiveWindow.ScrollRow = 15
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

# 3. Replace the `=FALSE()` formulas in H2:H49 with the literal boolean
#    value FALSE, keeping the existing "TRUE"/"FALSE" cell formatting.
$ws.Range("H2:H49").Value = $false
